# Auto-generated Excel COM-interop script to update the cryptos list
# (mirrors the data refresh captured in the commit's OOXML diff).
# Price/Volume columns (D/E) are stored as plain text in the sheet, even
# though many look numeric (thousands use "." as a separator, e.g.
# "96.604.68", and the cells carry no numeric format). A leading apostrophe
# is the standard Excel text-prefix marker: it forces the cell to remain
# text (keeping the "General" number format untouched) without the
# apostrophe itself becoming part of the stored value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'96.604.68"
$ws.Range("E2").Value = "'  +0.20%  "

$ws.Range("D3").Value = "'3.622.14"
$ws.Range("E3").Value = "'  +1.10%  "

$ws.Range("E4").Value = "'  -0.06%  "

$ws.Range("D5").Value = "'241.50"
$ws.Range("E5").Value = "'  +0.69%  "

$ws.Range("D6").Value = "'1.81"
$ws.Range("E6").Value = "'  +15.75%  "

$ws.Range("D7").Value = "'649.97"
$ws.Range("E7").Value = "'  -1.51%  "

$ws.Range("D8").Value = "'0.418"
$ws.Range("E8").Value = "'  +3.25%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'1.06"
$ws.Range("E9").Value = "'  +1.22%  "

$ws.Range("B10").Value = "USDC"
$ws.Range("C10").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D10").Value = "'0.999"
$ws.Range("E10").Value = "'  -0.05%  "

$ws.Range("D11").Value = "'3.613.98"
$ws.Range("E11").Value = "'  +0.96%  "

$ws.Range("D12").Value = "'44.15"
$ws.Range("E12").Value = "'  +1.49%  "

$ws.Range("E13").Value = "'  +0.08%  "

$ws.Range("D14").Value = "'6.51"
$ws.Range("E14").Value = "'  +1.73%  "

$ws.Range("D15").Value = "'4.289.13"
$ws.Range("E15").Value = "'  +0.84%  "

$ws.Range("D16").Value = "'96.352.68"
$ws.Range("E16").Value = "'  +0.09%  "

$ws.Range("D17").Value = "'0.0000258"
$ws.Range("E17").Value = "'  +0.32%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'8.69"
$ws.Range("E18").Value = "'  +11.78%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "'3.617.35"
$ws.Range("E19").Value = "'  +0.85%  "

$ws.Range("D20").Value = "'12.92"
$ws.Range("E20").Value = "'  +1.67%  "

$ws.Range("D21").Value = "'18.25"
$ws.Range("E21").Value = "'  +1.55%  "

$ws.Range("D22").Value = "'0.529"
$ws.Range("E22").Value = "'  +6.31%  "

$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "'3.43"
$ws.Range("E23").Value = "'  -0.53%  "

$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'508.67"
$ws.Range("E24").Value = "'  -0.72%  "

$ws.Range("D25").Value = "'0.0000204"
$ws.Range("E25").Value = "'  +2.20%  "

$ws.Range("D26").Value = "'6.83"
$ws.Range("E26").Value = "'  -1.00%  "

$ws.Range("D27").Value = "'101.55"
$ws.Range("E27").Value = "'  +4.59%  "

$ws.Range("D28").Value = "'13.20"
$ws.Range("E28").Value = "'  +2.86%  "

$ws.Range("D29").Value = "'0.169"
$ws.Range("E29").Value = "'  +13.67%  "

$ws.Range("E30").Value = "'  -2.03%  "

$ws.Range("D31").Value = "'11.97"
$ws.Range("E31").Value = "'  +3.12%  "

$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "'  -0.05%  "

$ws.Range("D33").Value = "'0.184"
$ws.Range("E33").Value = "'  +1.18%  "

$ws.Range("D34").Value = "'0.996"
$ws.Range("E34").Value = "'  -0.29%  "

$ws.Range("D35").Value = "'32.25"
$ws.Range("E35").Value = "'  +1.62%  "

$ws.Range("E36").Value = "'  +7.55%  "

$ws.Range("D37").Value = "'0.575"
$ws.Range("E37").Value = "'  +0.34%  "

$ws.Range("D38").Value = "'8.81"
$ws.Range("E38").Value = "'  +4.25%  "

$ws.Range("D39").Value = "'614.63"
$ws.Range("E39").Value = "'  +3.91%  "

$ws.Range("D40").Value = "'0.155"

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'0.947"
$ws.Range("E41").Value = "'  +4.06%  "

$ws.Range("B42").Value = "ImmutableX"
$ws.Range("C42").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D42").Value = "'1.91"
$ws.Range("E42").Value = "'  +2.83%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "'  -0.01%  "

$ws.Range("D44").Value = "'6.23"
$ws.Range("E44").Value = "'  +7.34%  "

$ws.Range("D45").Value = "'0.0444"
$ws.Range("E45").Value = "'  +5.34%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'35.21"
$ws.Range("E46").Value = "'  +2.64%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.418"
$ws.Range("E47").Value = "'  +17.20%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'2.28"
$ws.Range("E48").Value = "'  -0.39%  "

$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "'23.58"
$ws.Range("E49").Value = "'  +0.30%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'8.58"
$ws.Range("E50").Value = "'  +3.27%  "

$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").Value = "'54.04"
$ws.Range("E51").Value = "'  +0.42%  "
